$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 blank rows before the old "legend" footer block (old rows 31-35) ---
# This pushes the footer rows down to 34-38, leaving rows 30-33 empty except for
# the new Milestone 3 row we add at row 30.
$ws.Rows("31:33").Insert()

# --- Update existing row 29: change it from "Milestone 1 Testing & Evaluation"
#     to "Milestone 2 Testing & Evaluation" (+ matching rich-text description) ---
$ws.Range("C29").Value2 = "Milestone 2 Testing & Evaluation"

$d29Text = "Evaluation for milestone 2 components (Peak Load, Stats Priority, Optimizer Stats Generation Time Estimate models)"
$ws.Range("D29").Value2 = $d29Text
$d29 = $ws.Range("D29")
# Bold the "Peak Load, Stats Priority, Optimizer Stats Generation Time Estimate" run
$d29.Characters(40, 67).Font.Bold = $true
# Re-assert the trailing run's font so it gets its own explicit run properties
# (matches the formatting pattern used by the other rich-text cells in this sheet)
$tail = $d29.Characters(107, 8)
$tail.Font.Size = 11
$tail.Font.Name = "Calibri"

# --- Add new row 30: Milestone 3 / Write-Up / Compile Write Up ---
# Copy formatting from row 29 first so the new row's cells share the same style.
$ws.Range("A29:E29").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A30").Value2 = 29
$ws.Range("B30").Value2 = "Milestone 3"
$ws.Range("C30").Value2 = "Write-Up"
$ws.Range("D30").Value2 = "Compile Write Up"
$ws.Range("E30").Value2 = "N/A"

# --- Reset the view: scroll back to top and select D1 ---
[void]$ws.Range("D1").Select()
